$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Add the new time-registration entry in row 6 (task "UTD0104", role "Implenter")
$ws.Range("A6").Value = "UTD0104"
$ws.Range("B6").Value = "Implenter"
$ws.Range("C6").Value = 43893
$ws.Range("D6").Value = 0.60416666666666663
$ws.Range("E6").Value = 0.65625

# Move the selection to F4, matching the author's last cursor position
$ws.Range("F4").Select()
